$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.Formula = "=""" + $text + """"
    $range.Copy()
    $range.PasteSpecial(-4163)
}

Set-TextValue $ws.Range("D2") "29.449.64"
Set-TextValue $ws.Range("E2") "  +0.24%  "
Set-TextValue $ws.Range("D3") "1.869.86"
Set-TextValue $ws.Range("E3") "  -0.67%  "
Set-TextValue $ws.Range("D4") "1.001"
Set-TextValue $ws.Range("D5") "243.69"
Set-TextValue $ws.Range("E5") "  +0.26%  "
Set-TextValue $ws.Range("D6") "0.7052"
Set-TextValue $ws.Range("E6") "  -1.00%  "
Set-TextValue $ws.Range("E7") "  -0.12%  "
Set-TextValue $ws.Range("D8") "0.3148"
Set-TextValue $ws.Range("E8") "  -0.26%  "
Set-TextValue $ws.Range("D9") "0.07862"
Set-TextValue $ws.Range("E9") "  -2.04%  "
Set-TextValue $ws.Range("D10") "24.55"
Set-TextValue $ws.Range("E10") "  -2.17%  "
Set-TextValue $ws.Range("D11") "0.08022"
Set-TextValue $ws.Range("E11") "  -3.78%  "
Set-TextValue $ws.Range("D12") "1.891.26"
Set-TextValue $ws.Range("E12") "  -0.31%  "
Set-TextValue $ws.Range("D13") "5.206"
Set-TextValue $ws.Range("E13") "  -1.23%  "
Set-TextValue $ws.Range("D14") "93.96"
Set-TextValue $ws.Range("E14") "  -1.16%  "
Set-TextValue $ws.Range("D15") "0.7025"
Set-TextValue $ws.Range("E15") "  -2.26%  "
Set-TextValue $ws.Range("D16") "6.474"
Set-TextValue $ws.Range("E16") "  +1.63%  "
Set-TextValue $ws.Range("D17") "29.506.91"
Set-TextValue $ws.Range("E17") "  +0.39%  "
Set-TextValue $ws.Range("D18") "0.000008335"
Set-TextValue $ws.Range("E18") "  -3.80%  "
Set-TextValue $ws.Range("D19") "256.33"
Set-TextValue $ws.Range("E19") "  +5.48%  "
Set-TextValue $ws.Range("D20") "2.139.27"
Set-TextValue $ws.Range("E20") "  -0.54%  "
Set-TextValue $ws.Range("D21") "13.16"
Set-TextValue $ws.Range("E21") "  -1.39%  "
Set-TextValue $ws.Range("E22") "  -0.13%  "
Set-TextValue $ws.Range("D23") "7.614"
Set-TextValue $ws.Range("E23") "  -2.91%  "
Set-TextValue $ws.Range("D24") "1.001"
Set-TextValue $ws.Range("E24") "  -0.18%  "
Set-TextValue $ws.Range("D25") "0.1557"
Set-TextValue $ws.Range("E25") "  -1.23%  "
Set-TextValue $ws.Range("D26") "9.046"
Set-TextValue $ws.Range("E26") "  -0.55%  "
Set-TextValue $ws.Range("D27") "161.08"
Set-TextValue $ws.Range("E27") "  -1.44%  "
Set-TextValue $ws.Range("E28") "  +0.97%  "
Set-TextValue $ws.Range("D29") "1.502"
Set-TextValue $ws.Range("E29") "  -0.63%  "
Set-TextValue $ws.Range("D30") "4.330"
Set-TextValue $ws.Range("E30") "  -2.61%  "
Set-TextValue $ws.Range("D31") "4.264"
Set-TextValue $ws.Range("E31") "  -2.06%  "
Set-TextValue $ws.Range("D32") "1.203"
Set-TextValue $ws.Range("E32") "  -0.36%  "
Set-TextValue $ws.Range("D33") "0.05317"
Set-TextValue $ws.Range("E33") "  -1.57%  "
Set-TextValue $ws.Range("D34") "1.893"
Set-TextValue $ws.Range("E34") "  -2.85%  "
Set-TextValue $ws.Range("D35") "0.7468"
Set-TextValue $ws.Range("E35") "  -3.61%  "
Set-TextValue $ws.Range("D36") "1.166"
Set-TextValue $ws.Range("E36") "  -1.63%  "
Set-TextValue $ws.Range("D37") "2.721"
Set-TextValue $ws.Range("E37") "  +1.18%  "
Set-TextValue $ws.Range("D38") "0.01874"
Set-TextValue $ws.Range("E38") "  -0.92%  "
Set-TextValue $ws.Range("D39") "1.265.45"
Set-TextValue $ws.Range("E39") "  -0.66%  "
Set-TextValue $ws.Range("D40") "2.748"
Set-TextValue $ws.Range("E40") "  +0.10%  "
Set-TextValue $ws.Range("D41") "0.8984"
Set-TextValue $ws.Range("E41") "  -2.41%  "
Set-TextValue $ws.Range("D42") "108.78"
Set-TextValue $ws.Range("E42") "  -3.97%  "
Set-TextValue $ws.Range("D43") "5.997"
Set-TextValue $ws.Range("E43") "  -8.12%  "
Set-TextValue $ws.Range("D44") "71.45"
Set-TextValue $ws.Range("E44") "  -4.14%  "
Set-TextValue $ws.Range("D45") "1.001"
Set-TextValue $ws.Range("E45") "  -0.15%  "
Set-TextValue $ws.Range("D46") "0.00000000129"
Set-TextValue $ws.Range("E46") "  +0.50%  "
Set-TextValue $ws.Range("D47") "2.037.58"
Set-TextValue $ws.Range("E47") "  -0.26%  "
Set-TextValue $ws.Range("D48") "1.801"
Set-TextValue $ws.Range("E48") "  -0.83%  "
Set-TextValue $ws.Range("E49") "  -0.61%  "
Set-TextValue $ws.Range("D50") "9.503"
Set-TextValue $ws.Range("E50") "  -0.84%  "
Set-TextValue $ws.Range("D51") "0.4314"
Set-TextValue $ws.Range("E51") "  -1.57%  "
